$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H76").Value = 2418667.8
$ws.Range("I76").Value = 3341.1765
$ws.Range("K76").Value = 3341.1765
$ws.Range("M76").Value = -3026.1765

$ws.Range("H79").Value = 2418667.8
$ws.Range("I79").Value = 3341.1765
$ws.Range("K79").Value = 3341.1765
$ws.Range("M79").Value = -2249.1765

$ws.Range("H100").Value = 2161.5454
$ws.Range("I100").Value = 1868.1428
$ws.Range("J100").Value = 2675
$ws.Range("K100").Value = 1868.1428
$ws.Range("L100").Value = 2675
$ws.Range("M100").Value = -1327.1428
$ws.Range("N100").Value = -3757

$ws.Range("H137").Value = 1852.9412
$ws.Range("I137").Value = 2238.889
$ws.Range("J137").Value = 1418.75
$ws.Range("K137").Value = 6716.667
$ws.Range("L137").Value = 4256.25
$ws.Range("M137").Value = -4166.667
$ws.Range("N137").Value = -9356.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6563.727
$ws.Range("I32").Value = 5267.35
$ws.Range("J32").Value = 11139.177
$ws.Range("K32").Value = 5267.35
$ws.Range("L32").Value = 11139.177
$ws.Range("M32").Value = -4980.35
$ws.Range("N32").Value = -11713.177

$ws.Range("H45").Value = 3276.8147
$ws.Range("I45").Value = 2786.4
$ws.Range("K45").Value = 2786.4
$ws.Range("M45").Value = -2409.4

$ws.Range("H61").Value = 3356.9644
$ws.Range("I61").Value = 3225.8696
$ws.Range("J61").Value = 3960
$ws.Range("K61").Value = 3225.8696
$ws.Range("L61").Value = 3960
$ws.Range("M61").Value = -3013.8696
$ws.Range("N61").Value = -4384

$ws.Range("H102").Value = 1430.4667
$ws.Range("I102").Value = 1372.7693
$ws.Range("K102").Value = 1372.7693
$ws.Range("M102").Value = 249.2307000000001

$ws.Range("H136").Value = 3356.9644
$ws.Range("I136").Value = 3225.8696
$ws.Range("J136").Value = 3960
$ws.Range("K136").Value = 9677.6088
$ws.Range("L136").Value = 11880
$ws.Range("M136").Value = -7127.6088
$ws.Range("N136").Value = -16980

$ws.Range("H137").Value = 47734.832
$ws.Range("J137").Value = 48681.8
$ws.Range("L137").Value = 48681.8
$ws.Range("N137").Value = -58881.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 561.8125
$ws.Range("I94").Value = 419.08334
$ws.Range("K94").Value = 419.08334
$ws.Range("M94").Value = 31.91665999999998

$ws.Range("H105").Value = 1516957.1
$ws.Range("I105").Value = 1649.0714
$ws.Range("K105").Value = 1649.0714
$ws.Range("M105").Value = 97.92859999999996

$ws.Range("H134").Value = 3985.1282
$ws.Range("I134").Value = 3664.2424
$ws.Range("K134").Value = 10992.7272
$ws.Range("M134").Value = -8457.727200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3485.4194
$ws.Range("I31").Value = 847.0625
$ws.Range("J31").Value = 6299.6665
$ws.Range("K31").Value = 847.0625
$ws.Range("L31").Value = 6299.6665
$ws.Range("M31").Value = -552.0625
$ws.Range("N31").Value = -6889.6665

$ws.Range("H34").Value = 3485.4194
$ws.Range("I34").Value = 847.0625
$ws.Range("J34").Value = 6299.6665
$ws.Range("K34").Value = 847.0625
$ws.Range("L34").Value = 6299.6665
$ws.Range("M34").Value = -645.0625
$ws.Range("N34").Value = -6703.6665

$ws.Range("H62").Value = 5189.9
$ws.Range("I62").Value = 5199.857
$ws.Range("K62").Value = 5199.857
$ws.Range("M62").Value = -4575.857

$ws.Range("H65").Value = 5189.9
$ws.Range("I65").Value = 5199.857
$ws.Range("K65").Value = 25999.285
$ws.Range("M65").Value = -22879.285

$ws.Range("H107").Value = 1032.3334
$ws.Range("I107").Value = 365.83334
$ws.Range("K107").Value = 365.83334
$ws.Range("M107").Value = 1554.16666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 739.86
$ws.Range("I131").Value = 615
$ws.Range("J131").Value = 742.40814
$ws.Range("K131").Value = 1845
$ws.Range("L131").Value = 2227.22442
$ws.Range("M131").Value = 3195
$ws.Range("N131").Value = -12307.22442

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4475342.5
$ws.Range("I70").Value = 4750
$ws.Range("J70").Value = 10436133
$ws.Range("K70").Value = 4750
$ws.Range("L70").Value = 10436133
$ws.Range("M70").Value = -4480
$ws.Range("N70").Value = -10436673

$ws.Range("H73").Value = 4475342.5
$ws.Range("I73").Value = 4750
$ws.Range("J73").Value = 10436133
$ws.Range("K73").Value = 4750
$ws.Range("L73").Value = 10436133
$ws.Range("M73").Value = -3814
$ws.Range("N73").Value = -10438005

$ws.Range("H80").Value = 3444.2693
$ws.Range("I80").Value = 3210
$ws.Range("J80").Value = 3590.6875
$ws.Range("K80").Value = 3210
$ws.Range("L80").Value = 3590.6875
$ws.Range("M80").Value = -2212
$ws.Range("N80").Value = -5586.6875

$ws.Range("H83").Value = 3444.2693
$ws.Range("I83").Value = 3210
$ws.Range("J83").Value = 3590.6875
$ws.Range("K83").Value = 16050
$ws.Range("L83").Value = 17953.4375
$ws.Range("M83").Value = -11058
$ws.Range("N83").Value = -27937.4375

$ws.Range("H107").Value = 426
$ws.Range("I107").Value = 328.5
$ws.Range("J107").Value = 621
$ws.Range("K107").Value = 328.5
$ws.Range("L107").Value = 621
$ws.Range("M107").Value = 1591.5
$ws.Range("N107").Value = -4461

$ws.Range("H113").Value = 4744.9355
$ws.Range("I113").Value = 5835.952
$ws.Range("J113").Value = 2453.8
$ws.Range("K113").Value = 5835.952
$ws.Range("L113").Value = 2453.8
$ws.Range("M113").Value = -3665.952
$ws.Range("N113").Value = -6793.8

$ws.Range("H139").Value = 27206.357
$ws.Range("J139").Value = 27206.357
$ws.Range("L139").Value = 27206.357
$ws.Range("N139").Value = -37486.357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2054.0908
$ws.Range("I68").Value = 900
$ws.Range("K68").Value = 900
$ws.Range("M68").Value = -151

$ws.Range("H71").Value = 2054.0908
$ws.Range("I71").Value = 900
$ws.Range("K71").Value = 4500
$ws.Range("M71").Value = -756

$ws.Range("H93").Value = 883.75
$ws.Range("I93").Value = 883.75
$ws.Range("K93").Value = 883.75
$ws.Range("M93").Value = 364.25

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
